$wb = $excel.ActiveWorkbook

# --- "Range Status" sheet ---
$wsRange = $wb.Worksheets.Item("Range Status")

# Clear the "Species (perc.)" column C for rows 2 and 3 (now blank)
$wsRange.Range("C2").ClearContents()
$wsRange.Range("C3").ClearContents()

# Rows 4-7: "Species (no.)" column B reset to 0, "Species (perc.)" column C cleared
$wsRange.Range("B4").Value = 0
$wsRange.Range("C4").ClearContents()

$wsRange.Range("B5").Value = 0
$wsRange.Range("C5").ClearContents()

$wsRange.Range("B6").Value = 0
$wsRange.Range("C6").ClearContents()

$wsRange.Range("B7").Value = 0
$wsRange.Range("C7").ClearContents()

# --- "Species qualification" sheet ---
$wsQual = $wb.Worksheets.Item("Species qualification")

# "Range Analysis" row: "Selected for analysis" count reset to 0
$wsQual.Range("B5").Value = 0

# --- "High Priority break-up" sheet ---
$wsBreak = $wb.Worksheets.Item("High Priority break-up")

# "IUCN" row: add new high species counts
$wsBreak.Range("D2").Value = 1
$wsBreak.Range("E2").Value = 100
